$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# --- Donor cells used to clone number-format styles without creating new style entries ---
# I14 -> style 14 (integer #,##0)
# K14 -> style 15 (percent #,##0.0)
# D15 -> style 13 text placeholder holding shared string "0"
# E15 -> style 13 text placeholder holding shared string "***.*"

$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("I14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1

$ws.Range("I15").Value = 4

$ws.Range("K15").Value = 100

$ws.Range("L15").Value = 100

$ws.Range("M15").Value = 100

$ws.Range("N15").Value = -60

$ws.Range("D15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("C16").PasteSpecial(-4163)

$ws.Range("D16").Value = 2

$ws.Range("E16").Value = -100

$ws.Range("F16").Value = 2

$ws.Range("H16").Value = -60

$ws.Range("J16").Value = 34

$ws.Range("K16").Value = -38.235294117647

$ws.Range("L16").Value = -56.25

$ws.Range("M16").Value = -58

$ws.Range("N16").Value = -92.5

$ws.Range("I14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 4

$ws.Range("D17").Value = 4

$ws.Range("E17").Value = 0

$ws.Range("F17").Value = 9

$ws.Range("G17").Value = 11

$ws.Range("H17").Value = -18.181818181818

$ws.Range("I17").Value = 55

$ws.Range("J17").Value = 53

$ws.Range("K17").Value = 3.773584905660

$ws.Range("L17").Value = -15.384615384615

$ws.Range("M17").Value = 450

$ws.Range("N17").Value = -23.611111111111

$ws.Range("C18").Value = 4

$ws.Range("D15").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D18").PasteSpecial(-4163)

$ws.Range("E15").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("F18").Value = 7

$ws.Range("G18").Value = 6

$ws.Range("H18").Value = 16.666666666666

$ws.Range("I18").Value = 38

$ws.Range("K18").Value = -15.555555555555

$ws.Range("L18").Value = -51.282051282051

$ws.Range("M18").Value = -38.709677419354

$ws.Range("N18").Value = -85.496183206106

$ws.Range("C19").Value = 6

$ws.Range("D19").Value = 10

$ws.Range("E19").Value = -40

$ws.Range("F19").Value = 27

$ws.Range("G19").Value = 36

$ws.Range("H19").Value = -25

$ws.Range("I19").Value = 127

$ws.Range("J19").Value = 160

$ws.Range("K19").Value = -20.625

$ws.Range("L19").Value = -34.196891191709

$ws.Range("M19").Value = -14.189189189189

$ws.Range("N19").Value = -13.605442176870

$ws.Range("I14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1

$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 0

$ws.Range("F20").Value = 5

$ws.Range("G20").Value = 12

$ws.Range("H20").Value = -58.333333333333

$ws.Range("I20").Value = 15

$ws.Range("J20").Value = 51

$ws.Range("K20").Value = -70.588235294117

$ws.Range("L20").Value = -57.142857142857

$ws.Range("M20").Value = -40

$ws.Range("N20").Value = -96.25

$ws.Range("C21").Value = 16

$ws.Range("D21").Value = 17

$ws.Range("E21").Value = -5.882352941176

$ws.Range("F21").Value = 51

$ws.Range("G21").Value = 70

$ws.Range("H21").Value = -27.142857142857

$ws.Range("I21").Value = 261

$ws.Range("J21").Value = 347

$ws.Range("K21").Value = -24.783861671469

$ws.Range("L21").Value = -38.151658767772

$ws.Range("M21").Value = -12.121212121212

$ws.Range("N21").Value = -77.749360613810

$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 3

$ws.Range("F22").Value = 3

$ws.Range("D15").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("G22").PasteSpecial(-4163)

$ws.Range("E15").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("H22").PasteSpecial(-4163)

$ws.Range("I22").Value = 11

$ws.Range("K22").Value = 22.222222222222

$ws.Range("L22").Value = 37.5

$ws.Range("M22").Value = 10

$ws.Range("D15").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("F23").PasteSpecial(-4163)

$ws.Range("H23").Value = -100

$ws.Range("C24").Value = 26

$ws.Range("D24").Value = 23

$ws.Range("E24").Value = 13.043478260869

$ws.Range("F24").Value = 107

$ws.Range("G24").Value = 89

$ws.Range("H24").Value = 20.224719101123

$ws.Range("I24").Value = 483

$ws.Range("J24").Value = 378

$ws.Range("K24").Value = 27.777777777777

$ws.Range("L24").Value = -16.291161178509

$ws.Range("M24").Value = 71.276595744680

$ws.Range("C25").Value = 11

$ws.Range("D25").Value = 15

$ws.Range("E25").Value = -26.666666666666

$ws.Range("F25").Value = 64

$ws.Range("G25").Value = 59

$ws.Range("H25").Value = 8.474576271186

$ws.Range("I25").Value = 325

$ws.Range("J25").Value = 254

$ws.Range("K25").Value = 27.952755905511

$ws.Range("L25").Value = -28.256070640176

$ws.Range("C26").Value = 3

$ws.Range("E26").Value = -25

$ws.Range("F26").Value = 17

$ws.Range("G26").Value = 18

$ws.Range("H26").Value = -5.555555555555

$ws.Range("I26").Value = 82

$ws.Range("J26").Value = 75

$ws.Range("K26").Value = 9.333333333333

$ws.Range("L26").Value = -21.904761904761

$ws.Range("M26").Value = 18.840579710144

$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

$ws.Range("I14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1

$ws.Range("D15").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("G27").PasteSpecial(-4163)

$ws.Range("E15").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("H27").PasteSpecial(-4163)

$ws.Range("I27").Value = 6

$ws.Range("K27").Value = 0

$ws.Range("L27").Value = 20

$ws.Range("C28").Value = 1

$ws.Range("E28").Value = -50

$ws.Range("F28").Value = 7

$ws.Range("G28").Value = 10

$ws.Range("H28").Value = -30

$ws.Range("I28").Value = 18

$ws.Range("J28").Value = 34

$ws.Range("K28").Value = -47.058823529411

$ws.Range("L28").Value = 28.571428571428

$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$ws.Range("I14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1

$ws.Range("D15").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("G29").PasteSpecial(-4163)

$ws.Range("E15").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("H29").PasteSpecial(-4163)

$ws.Range("I14").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("I29").Value = 1

$ws.Range("K29").Value = -66.666666666666

$ws.Range("L29").Value = 0

$ws.Range("N29").Value = -75

$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

$ws.Range("I14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1

$ws.Range("D15").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("G30").PasteSpecial(-4163)

$ws.Range("E15").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("H30").PasteSpecial(-4163)

$ws.Range("I14").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("I30").Value = 1

$ws.Range("K30").Value = -66.666666666666

$ws.Range("L30").Value = 0

$ws.Range("N30").Value = -75
